# Insert a new row at position 12 (shifts existing rows 12-105 down to 13-106)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(12).Insert()

# Populate the newly inserted row 12 with its data
$ws.Range("A12").Value = 3
$ws.Range("B12").Value = 'Femacal de La Calera'
$ws.Range("C12").Value = 'Coquimbo'
$ws.Range("D12").Value = 44532
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 100112026
$ws.Range("G12").Value = 'Haba'
$ws.Range("H12").Value = 'Sin especificar'
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 88
$ws.Range("K12").Value = 8000
$ws.Range("L12").Value = 8500
$ws.Range("M12").Value = 8273
$ws.Range("N12").Value = '$/malla 25 kilos'
$ws.Range("O12").Value = 'Provincia de Limarí'
$ws.Range("P12").Value = 331
$ws.Range("Q12").Value = 25
$ws.Range("R12").Value = 'Hortaliza'
